$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 changes from "Document StateEngine" to "Document TileMap"
$ws.Range("C2").Value = "Document TileMap"

# C3 changes from "Document TileMap" to "Document SDLUtilityTool"
$ws.Range("C3").Value = "Document SDLUtilityTool"

# B2 is new: "Document the TileMap, b2world, axis orientations"
$ws.Range("B2").Value = "Document the TileMap, b2world, axis orientations"

# Page setup: explicit portrait orientation (adds <pageSetup orientation="portrait"/>)
$ws.PageSetup.Orientation = 1

# Update selection to C2:C3 with active cell C2
$ws.Range("C2:C3").Select()
